$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh capital-structure / financial metrics for the India Banks (Regional) rows

# Row 2
$ws.Range("D2").Value = 0.2285
$ws.Range("F2").Value = 0.164
$ws.Range("I2").Value = 0.003397303966033762
$ws.Range("J2").Value = 0.002788745102350375
$ws.Range("K2").Value = 544.28
$ws.Range("L2").Value = 0.3478271983640082
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 1349.6
$ws.Range("V2").Value = 0.09601183785552693
$ws.Range("W2").Value = 0.163594470046083
$ws.Range("X2").Value = 0.05849451179760351
$ws.Range("Y2").Value = 0.1050999582484795
$ws.Range("Z2").Value = 0.2461046990117819
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04998027871390719
$ws.Range("AC2").Value = -0.04967108931316869
$ws.Range("AD2").Value = 4136.3
$ws.Range("AE2").Value = 149.9694937697519
$ws.Range("AF2").Value = 4286.269493769752
$ws.Range("AG2").Value = 2936.669493769752
$ws.Range("AH2").Value = 0.233674970823164
$ws.Range("AI2").Value = 0.5307239310273503
$ws.Range("AJ2").Value = 0.1728136833730804
$ws.Range("AK2").Value = 0.4365710990393831
$ws.Range("AN2").Value = 117.1424525630133
$ws.Range("AP2").Value = 83.16820996232657

# Row 3
$ws.Range("D3").Value = 0.297
$ws.Range("I3").Value = 0.0154063418587437
$ws.Range("J3").Value = 0.0154063418587437
$ws.Range("K3").Value = 5.98
$ws.Range("L3").Value = 0.1196
$ws.Range("U3").Value = 146.7
$ws.Range("V3").Value = 2.916500994035785
$ws.Range("W3").Value = 0.05259454705364996
$ws.Range("X3").Value = 0.05849451179760351
$ws.Range("Y3").Value = -0.005899964743953552
$ws.Range("Z3").Value = 1.519830080149565
$ws.Range("AA3").Value = 0.02341502178198603
$ws.Range("AB3").Value = 0.04998027871390719
$ws.Range("AC3").Value = -0.02656525693192116
$ws.Range("AD3").Value = 20.4
$ws.Range("AE3").Value = 2.998414535314076
$ws.Range("AF3").Value = 23.39841453531407
$ws.Range("AG3").Value = -123.3015854646859
$ws.Range("AH3").Value = 0.3174887096669122
$ws.Range("AI3").Value = 0.1688216608664906
$ws.Range("AJ3").Value = 1.689026131142485
$ws.Range("AK3").Value = 15.21943896070054
$ws.Range("AN3").Value = 14.89051094890511
$ws.Range("AP3").Value = -90.00115727349336

# Row 4
$ws.Range("F4").Value = 0.24
$ws.Range("I4").Value = 0.007781787248309645
$ws.Range("J4").Value = 0.005825700966123649
$ws.Range("K4").Value = 369.7
$ws.Range("L4").Value = 0.437359517331125
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 425.4
$ws.Range("V4").Value = 0.04795616981940342
$ws.Range("W4").Value = 0.1863595120475854
$ws.Range("X4").Value = 0.05211368738715745
$ws.Range("Y4").Value = 0.134245824660428
$ws.Range("Z4").Value = 0.2242168954651075
$ws.Range("AA4").Value = 0.001306220584532322
$ws.Range("AB4").Value = 0.04843462425646064
$ws.Range("AC4").Value = -0.04712840367192832
$ws.Range("AD4").Value = 1905.4
$ws.Range("AE4").Value = 42.11027619501929
$ws.Range("AF4").Value = 1947.510276195019
$ws.Range("AG4").Value = 1522.110276195019
$ws.Range("AH4").Value = 0.1800231488192967
$ws.Range("AI4").Value = 0.4621183891456067
$ws.Range("AJ4").Value = 0.1464594158543496
$ws.Range("AK4").Value = 0.4017277172695642
$ws.Range("AN4").Value = 127.0266666666667
$ws.Range("AP4").Value = 101.4740184130013

# Row 5
$ws.Range("B5").Value = "AU Small Finance Bank Limited (BSE:540611)"
$ws.Range("F5").Value = 0.164
$ws.Range("I5").Value = -0.000355746682121268
$ws.Range("J5").Value = -0.0002674574658436536
$ws.Range("K5").Value = 113.6
$ws.Range("L5").Value = 0.3256880733944954
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("U5").Value = 588.3
$ws.Range("V5").Value = 0.1643341992793095
$ws.Range("W5").Value = 0.2282958199356913
$ws.Range("X5").Value = 0.05674558327487197
$ws.Range("Y5").Value = 0.1715502366608193
$ws.Range("Z5").Value = 0.2259110634385851
$ws.Range("AA5").Value = -0.00006042160053332882
$ws.Range("AB5").Value = 0.04961066771263536
$ws.Range("AC5").Value = -0.04967108931316869
$ws.Range("AD5").Value = 1376.7
$ws.Range("AE5").Value = 47.57042221361949
$ws.Range("AF5").Value = 1424.270422213619
$ws.Range("AG5").Value = 835.9704222136195
$ws.Range("AH5").Value = 0.2846166900893808
$ws.Range("AI5").Value = 0.6805344502442608
$ws.Range("AJ5").Value = 0.1893104512325247
$ws.Range("AK5").Value = 0.5556206674485111
$ws.Range("AN5").Value = 146.6134185303514
$ws.Range("AP5").Value = 89.02773399506064

# Row 6
$ws.Range("B6").Value = "CSB Bank Limited (BSE:542867)"
$ws.Range("D6").Value = 0.16
$ws.Range("K6").Value = 12.4
$ws.Range("L6").Value = 0.1371681415929203
$ws.Range("U6").Value = 86.5
$ws.Range("V6").Value = 0.1669561860644663
$ws.Range("W6").Value = 0.05703771849126035
$ws.Range("X6").Value = 0.06128686938261257
$ws.Range("Y6").Value = -0.004249150891352221
$ws.Range("Z6").Value = 0.4657393096342092
$ws.Range("AB6").Value = 0.05050480544957378
$ws.Range("AC6").Value = -0.05050480544957378
$ws.Range("AD6").Value = 296.7
$ws.Range("AF6").Value = 296.7
$ws.Range("AG6").Value = 210.2
$ws.Range("AH6").Value = 0.3641384388807069
$ws.Range("AI6").Value = 0.511551724137931
$ws.Range("AJ6").Value = 0.2886173280241658
$ws.Range("AK6").Value = 0.4259371833839919

# Row 7
$ws.Range("B7").Value = "Ujjivan Small Finance Bank Limited (BSE:542904)"
$ws.Range("F7").Value = -0.298
$ws.Range("I7").Value = -0.008285176574727746
$ws.Range("J7").Value = -0.006280222812871922
$ws.Range("K7").Value = 42.6
$ws.Range("L7").Value = 0.1849761181068172
$ws.Range("M7").Value = -0
$ws.Range("N7").Value = -0
$ws.Range("O7").Value = -0
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = -0
$ws.Range("U7").Value = 102.7
$ws.Range("V7").Value = 0.09896887347017443
$ws.Range("W7").Value = 0.163594470046083
$ws.Range("X7").Value = 0.06129015821006738
$ws.Range("Y7").Value = 0.1023043118360156
$ws.Range("Z7").Value = 0.2817847920433156
$ws.Range("AA7").Value = -0.001769671279310801
$ws.Range("AB7").Value = 0.05050538096299727
$ws.Range("AC7").Value = -0.05227505224230807
$ws.Range("AD7").Value = 537.1
$ws.Range("AE7").Value = 57.29038082579901
$ws.Range("AF7").Value = 594.3903808257991
$ws.Range("AG7").Value = 491.6903808257991
$ws.Range("AH7").Value = 0.3641896232027614
$ws.Range("AI7").Value = 0.5658218215749333
$ws.Range("AJ7").Value = 0.321494359445565
$ws.Range("AK7").Value = 0.5187754494800789
$ws.Range("AN7").Value = 56.24083769633508
$ws.Range("AP7").Value = 51.48590375139257

# buybacks_cash_returned (T) no longer populated for these rows in the refreshed export
$ws.Range("T2").Value = $null
$ws.Range("T4").Value = $null
$ws.Range("T5").Value = $null
$ws.Range("T7").Value = $null
